$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.106.69'
$ws.Range("E2").Value = '  +0.62%  '
$ws.Range("D3").Value = '3.337.26'
$ws.Range("E3").Value = '  +0.46%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '584.41'
$ws.Range("E5").Value = '  +0.61%  '
$ws.Range("D6").Value = '176.88'
$ws.Range("E6").Value = '  +1.83%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  +1.41%  '
$ws.Range("E9").Value = '  +4.87%  '
$ws.Range("E10").Value = '  +1.48%  '
$ws.Range("D11").Value = '47.99'
$ws.Range("E11").Value = '  +5.75%  '
$ws.Range("D13").Value = '693.51'
$ws.Range("E13").Value = '  +3.98%  '
$ws.Range("D14").Value = '3.879.23'
$ws.Range("E14").Value = '  +0.50%  '
$ws.Range("D15").Value = '8.44'
$ws.Range("E15").Value = '  +0.81%  '
$ws.Range("D16").Value = '68.120.19'
$ws.Range("E16").Value = '  +0.39%  '
$ws.Range("E17").Value = '  +1.35%  '
$ws.Range("D18").Value = '3.289.62'
$ws.Range("E18").Value = '  -0.89%  '
$ws.Range("E19").Value = '  +0.53%  '
$ws.Range("D20").Value = '11.18'
$ws.Range("E20").Value = '  +2.78%  '
$ws.Range("D21").Value = '0.896'
$ws.Range("E21").Value = '  +0.99%  '
$ws.Range("E22").Value = '  +1.02%  '
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("D24").Value = '100.71'
$ws.Range("E24").Value = '  +3.29%  '
$ws.Range("D25").Value = '3.91'
$ws.Range("E25").Value = '  +1.99%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.70'
$ws.Range("E26").Value = '  +1.08%  '
$ws.Range("E27").Value = '  +2.29%  '
$ws.Range("E28").Value = '  -0.82%  '
$ws.Range("E29").Value = '  +1.77%  '
$ws.Range("D30").Value = '6.96'
$ws.Range("E30").Value = '  -4.72%  '
$ws.Range("D31").Value = '567.54'
$ws.Range("E31").Value = '  -3.07%  '
$ws.Range("D32").Value = '11.06'
$ws.Range("E32").Value = '  +1.19%  '
$ws.Range("E33").Value = '  +1.83%  '
$ws.Range("D34").Value = '3.738.10'
$ws.Range("E34").Value = '  +0.69%  '
$ws.Range("D35").Value = '57.46'
$ws.Range("E36").Value = '  -0.09%  '
$ws.Range("E37").Value = '  +2.68%  '
$ws.Range("E38").Value = '  +3.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.20'
$ws.Range("E39").Value = '  +7.69%  '
$ws.Range("D40").Value = '3.18'
$ws.Range("E40").Value = '  +2.61%  '
$ws.Range("E41").Value = '  +0.22%  '
$ws.Range("D42").Value = '0.0₃0674'
$ws.Range("E42").Value = '  +2.01%  '
$ws.Range("E43").Value = '  +1.06%  '
$ws.Range("E44").Value = '  +0.65%  '
$ws.Range("E45").Value = '  +1.33%  '
$ws.Range("E46").Value = '  +3.08%  '
$ws.Range("E47").Value = '  +1.27%  '
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("E49").Value = '  -0.52%  '
$ws.Range("D50").Value = '130.87'
$ws.Range("E50").Value = '  +3.27%  '
$ws.Range("E51").Value = '  +0.90%  '
